$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.937.08"
$ws.Range("E2").Value = "  +5.85%  "
$ws.Range("D3").Value = "2.508.96"
$ws.Range("E3").Value = "  +6.46%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "490.16"
$ws.Range("E5").Value = "  +7.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.66"
$ws.Range("E6").Value = "  +13.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.514"
$ws.Range("E8").Value = "  +8.23%  "
$ws.Range("D9").Value = "2.529.92"
$ws.Range("E9").Value = "  +7.29%  "
$ws.Range("E10").Value = "  +7.21%  "
$ws.Range("E11").Value = "  +6.84%  "
$ws.Range("E12").Value = "  +7.02%  "
$ws.Range("E13").Value = "  +2.00%  "
$ws.Range("D14").Value = "2.945.47"
$ws.Range("E14").Value = "  +6.60%  "
$ws.Range("D15").Value = "56.015.18"
$ws.Range("E15").Value = "  +6.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.08"
$ws.Range("E16").Value = "  +9.78%  "
$ws.Range("E17").Value = "  +11.44%  "
$ws.Range("D18").Value = "2.526.99"
$ws.Range("E18").Value = "  +6.34%  "
$ws.Range("E19").Value = "  +7.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.26"
$ws.Range("E20").Value = "  +13.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.11"
$ws.Range("E21").Value = "  +6.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.85"
$ws.Range("E23").Value = "  +10.91%  "
$ws.Range("E24").Value = "  +6.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.169"
$ws.Range("E25").Value = "  +14.14%  "
$ws.Range("E26").Value = "  +9.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "2.595.36"
$ws.Range("E28").Value = "  +6.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.54"
$ws.Range("E29").Value = "  +7.07%  "
$ws.Range("E30").Value = "  +15.59%  "
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "149.26"
$ws.Range("E32").Value = "  +2.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.50"
$ws.Range("E33").Value = "  +6.96%  "
$ws.Range("E34").Value = "  +9.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.29"
$ws.Range("E35").Value = "  +7.68%  "
$ws.Range("E36").Value = "  +13.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.72"
$ws.Range("E37").Value = "  +9.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.873"
$ws.Range("E38").Value = "  +6.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.31"
$ws.Range("E39").Value = "  +4.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.53"
$ws.Range("E40").Value = "  +9.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.623"
$ws.Range("E41").Value = "  +5.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0558"
$ws.Range("E42").Value = "  +7.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.996"
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("E44").Value = "  +12.03%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "269.80"
$ws.Range("E45").Value = "  +31.71%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.86"
$ws.Range("E46").Value = "  +17.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.18"
$ws.Range("E47").Value = "  +0.87%  "
$ws.Range("E48").Value = "  +7.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0906"
$ws.Range("E49").Value = "  +7.11%  "
$ws.Range("D50").Value = "1.960.49"
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.81"
$ws.Range("E51").Value = "  +9.02%  "
